# "Generate Report for Handback"
#
# Before this edit the status rows for the two source files
# (308f3b28-...md and 6b083ed6-...md) were still flagged "Ready for
# handoff" and had no "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" recorded. This script fills in the
# handback report: flips the status to "Handed back: in sync with
# en-US", records the target (.md) and handback (.xlf) file names
# (with hyperlinks back to the source .md on GitHub for the target
# file), and stamps the handback datetime - once for zh-cn, once for
# de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

$file1 = "308f3b28-a034-47cb-9b43-3f727e5c105a.md"
$file2 = "6b083ed6-3a0c-4639-82f7-7375ea43d6fd.md"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b644666ca9740b3caff25b4f21c2a38a1e5ea95/e2e/$file1"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b644666ca9740b3caff25b4f21c2a38a1e5ea95/e2e/$file2"

$zhcnXlf1 = "308f3b28-a034-47cb-9b43-3f727e5c105a.b81cf6c620ca181fd01919b663a7522974338052.zh-cn.xlf"
$zhcnXlf2 = "6b083ed6-3a0c-4639-82f7-7375ea43d6fd.64c87aaf6e6358b20f8c461e3cfbac2ec25ab11a.zh-cn.xlf"
$dedeXlf1 = "308f3b28-a034-47cb-9b43-3f727e5c105a.b81cf6c620ca181fd01919b663a7522974338052.de-de.xlf"
$dedeXlf2 = "6b083ed6-3a0c-4639-82f7-7375ea43d6fd.64c87aaf6e6358b20f8c461e3cfbac2ec25ab11a.de-de.xlf"

$zhcnHandbackTime = "2016-09-05 03:09:32"
$dedeHandbackTime = "2016-09-05 03:09:40"

# ---- Overview sheet: status text flips for both locale columns ----
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ==================== zh-cn sheet ====================

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = $file1
$zhcn.Range("J2").Value = $zhcnXlf1
$zhcn.Range("K2").Value = $zhcnHandbackTime

$zhcn.Range("I3").Value = $file2
$zhcn.Range("J3").Value = $zhcnXlf2
$zhcn.Range("K3").Value = $zhcnHandbackTime

# Rebuild the hyperlinks collection so the new Latest Target File
# links land in the same row-major order as the existing source file
# links (A2, I2, A3, I3).
$zhcn.Range("A2").Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $file1)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $file1)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $file2)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $file2)

$zhcn.Range("I2").Style = "HyperLink"
$zhcn.Range("I3").Style = "HyperLink"

$zhcn.Columns.Item(3).ColumnWidth  = 29 + 1/6
$zhcn.Columns.Item(9).ColumnWidth  = 39 + 1/6
$zhcn.Columns.Item(10).ColumnWidth = 39 + 1/6

# ==================== de-de sheet ====================

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $file1
$dede.Range("J2").Value = $dedeXlf1
$dede.Range("K2").Value = $dedeHandbackTime

$dede.Range("I3").Value = $file2
$dede.Range("J3").Value = $dedeXlf2
$dede.Range("K3").Value = $dedeHandbackTime

$dede.Range("A2").Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $file1)
$dede.Hyperlinks.Add($dede.Range("I2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $file1)
$dede.Hyperlinks.Add($dede.Range("A3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $file2)
$dede.Hyperlinks.Add($dede.Range("I3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $file2)

$dede.Range("I2").Style = "HyperLink"
$dede.Range("I3").Style = "HyperLink"

$dede.Columns.Item(3).ColumnWidth  = 29 + 1/6
$dede.Columns.Item(9).ColumnWidth  = 39 + 1/6
$dede.Columns.Item(10).ColumnWidth = 39 + 1/6

# Overview "Status" columns widen the same way the locale sheets' do.
$overview.Columns.Item(5).ColumnWidth = 29 + 1/6
$overview.Columns.Item(6).ColumnWidth = 29 + 1/6
